$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "dafd"
$ws.Range("C2").Value = "d"
$ws.Range("E2").Value = 44594

# Delete row 3 entirely (shifts rows up, removing the row)
$ws.Rows(3).Delete()
